# Edit "Welcome Training.docx": restructure the training-session paragraph and
# trim the closing line, per the commit "edited welcome training doc".

$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# 1) A second space slipped in between "email." and "After" when the paragraph
#    was re-flowed.
$null = $find.Execute("via email. After that", $true, $false, $false, $false, $false, $true, 1, $false, "via email.  After that", 2)

# 2) The manual line break after "...30 min training session." is removed so the
#    sentence about the session rules runs on in the same paragraph.
$null = $find.Execute("training session.^l", $true, $false, $false, $false, $false, $true, 1, $false, "training session. ", 2)

# 3) Same for the line break after "...code coaching." before "Then I will schedule...".
$null = $find.Execute("code coaching.^l", $true, $false, $false, $false, $false, $true, 1, $false, "code coaching. ", 2)

# 4) Turn the line break before "Also, if you have a" into an actual paragraph
#    break, splitting the intro paragraph into two.
$null = $find.Execute("shadow me.^lAlso,", $true, $false, $false, $false, $false, $true, 1, $false, "shadow me.^l^pAlso,", 2)

# 5) The non-breaking space that used to sit between "a" and "preference" becomes
#    a normal space now that "Also, if you have a " starts its own paragraph.
$null = $find.Execute("have a preference", $true, $false, $false, $false, $false, $true, 1, $false, "have a preference", 2)

# 6) Drop the "Anyways, that it is! " lead-in from the final sentence.
$null = $find.Execute("Anyways, that it is! Hope", $true, $false, $false, $false, $false, $true, 1, $false, "Hope", 2)
